# Atualização automática de preços de eletricidade
# Updates row 2 of the SpotPTTable (daily hourly spot prices) with the
# new day's figures, per the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (serial date 45930 -> 45931, i.e. 2025-09-30 -> 2025-10-01)
$ws.Range("A2").Value = 45931

# Hourly spot prices, 0h-1h .. 23h-24h
$ws.Range("B2").Value = 103.41
$ws.Range("C2").Value = 105.63
$ws.Range("D2").Value = 101.55
$ws.Range("E2").Value = 97.63
$ws.Range("F2").Value = 101.45
$ws.Range("G2").Value = 101.59
$ws.Range("H2").Value = 103.48
$ws.Range("I2").Value = 113.42
$ws.Range("J2").Value = 123.29
$ws.Range("K2").Value = 95.91
$ws.Range("L2").Value = 67.26000000000001
$ws.Range("M2").Value = 41.03
$ws.Range("N2").Value = 23.92
$ws.Range("O2").Value = 16.28
$ws.Range("P2").Value = 14.02
$ws.Range("Q2").Value = 19.25
$ws.Range("R2").Value = 36.97
$ws.Range("S2").Value = 53.19
$ws.Range("T2").Value = 86.81
$ws.Range("U2").Value = 125.06
$ws.Range("V2").Value = 217.9
$ws.Range("W2").Value = 129.41
$ws.Range("X2").Value = 108.47
$ws.Range("Y2").Value = 103.35
$ws.Range("Z2").Value = 87.09

# Slot_4h_max label unchanged (20h-24h); Slot_4h_price updated
$ws.Range("AB2").Value = 139.78

# Slot_2h_frist label unchanged (20h-22h); Slot_2h_frist_price updated
$ws.Range("AD2").Value = 173.66

# Slot_2h_second label and price updated
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 109.6

# Slot_min_price label updated
$ws.Range("AG2").Value = "10h-18h"
